# "Aggregated variables but new delivery variables must be converted to TAF!"
#
# The "Init" sheet's Variables List Indices row (row 13) records the block
# bounds used to pull variables out of the trend-report listing. Its
# "Lower Right Cell" (D13) moves down one row, from E177 to E178, to make
# room for the newly aggregated delivery variables that still need a TAF
# conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Variables List Indices -> Lower Right Cell: E177 -> E178
$ws.Range("D13").Value = "E178"

# Leave the sheet scrolled/selected where the edit was made (column B in view,
# D13 the active cell), matching the author's saved view state.
$ws.Range("D13").Select()
$excel.ActiveWindow.ScrollColumn = 2
